# Refresh market-price-derived columns (H:N) for a set of Leve rows
# across multiple worksheets, per the scheduled-runner price update.
$wb = $excel.ActiveWorkbook

# ALC!A19 row
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4478.6
$ws.Range("I19").Value = 6299.8237
$ws.Range("J19").Value = 608.5
$ws.Range("K19").Value = 6299.8237
$ws.Range("L19").Value = 608.5
$ws.Range("M19").Value = -6124.8237
$ws.Range("N19").Value = -958.5

# ALC!A98 row
$ws.Range("H98").Value = 960.8823
$ws.Range("I98").Value = 918.7273
$ws.Range("J98").Value = 1038.1666
$ws.Range("K98").Value = 918.7273
$ws.Range("L98").Value = 1038.1666
$ws.Range("M98").Value = 579.2727
$ws.Range("N98").Value = -4034.1666

# ALC!A116 row
$ws.Range("H116").Value = 63771.766
$ws.Range("I116").Value = 71154.664
$ws.Range("J116").Value = 8400
$ws.Range("K116").Value = 71154.664
$ws.Range("L116").Value = 8400
$ws.Range("M116").Value = -67712.664
$ws.Range("N116").Value = -15284

# ALC!A122 row
$ws.Range("H122").Value = 960.8823
$ws.Range("I122").Value = 918.7273
$ws.Range("J122").Value = 1038.1666
$ws.Range("K122").Value = 2756.1819
$ws.Range("L122").Value = 3114.4998
$ws.Range("M122").Value = -306.1819
$ws.Range("N122").Value = -8014.4998

# ARM!A45 row
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1165.35
$ws.Range("I45").Value = 1009.75
$ws.Range("J45").Value = 1398.75
$ws.Range("K45").Value = 1009.75
$ws.Range("L45").Value = 1398.75
$ws.Range("M45").Value = -632.75
$ws.Range("N45").Value = -2152.75

# ARM!A132 row
$ws.Range("H132").Value = 2885.4312
$ws.Range("I132").Value = 2683.1843
$ws.Range("J132").Value = 3269.7
$ws.Range("K132").Value = 8049.5529
$ws.Range("L132").Value = 9809.099999999999
$ws.Range("M132").Value = -5519.5529
$ws.Range("N132").Value = -14869.1

# BSM!A80 row
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1034.5
$ws.Range("I80").Value = 300
$ws.Range("J80").Value = 1181.4
$ws.Range("K80").Value = 300
$ws.Range("L80").Value = 1181.4
$ws.Range("M80").Value = 698
$ws.Range("N80").Value = -3177.4

# BSM!A83 row
$ws.Range("H83").Value = 1034.5
$ws.Range("I83").Value = 300
$ws.Range("J83").Value = 1181.4
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 5907
$ws.Range("M83").Value = 3492
$ws.Range("N83").Value = -15891

# BSM!A134 row
$ws.Range("H134").Value = 20279.729
$ws.Range("I134").Value = 23861.936
$ws.Range("J134").Value = 7604.231
$ws.Range("K134").Value = 71585.808
$ws.Range("L134").Value = 22812.693
$ws.Range("M134").Value = -69050.808
$ws.Range("N134").Value = -27882.693

# CRP!A16 row
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 718.2
$ws.Range("I16").Value = 718.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 718.2
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("M16").Value = -431.2

# CRP!A31 row
$ws.Range("H31").Value = 3186.1345
$ws.Range("I31").Value = 2198.8333
$ws.Range("J31").Value = 4532.4546
$ws.Range("K31").Value = 2198.8333
$ws.Range("L31").Value = 4532.4546
$ws.Range("M31").Value = -1903.8333
$ws.Range("N31").Value = -5122.4546

# CRP!A34 row
$ws.Range("H34").Value = 3186.1345
$ws.Range("I34").Value = 2198.8333
$ws.Range("J34").Value = 4532.4546
$ws.Range("K34").Value = 2198.8333
$ws.Range("L34").Value = 4532.4546
$ws.Range("M34").Value = -1996.8333
$ws.Range("N34").Value = -4936.4546

# CRP!A94 row
$ws.Range("H94").Value = 7135.7334
$ws.Range("I94").Value = 2405.1428
$ws.Range("J94").Value = 11275
$ws.Range("K94").Value = 2405.1428
$ws.Range("L94").Value = 11275
$ws.Range("M94").Value = -1954.1428
$ws.Range("N94").Value = -12177

# CRP!A99 row
$ws.Range("H99").Value = 44016.582
$ws.Range("I99").Value = 73209.57000000001
$ws.Range("J99").Value = 3146.4
$ws.Range("K99").Value = 73209.57000000001
$ws.Range("L99").Value = 3146.4
$ws.Range("M99").Value = -71711.57000000001
$ws.Range("N99").Value = -6142.4

# CRP!A113 row
$ws.Range("H113").Value = 718.2
$ws.Range("I113").Value = 718.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 718.2
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("M113").Value = 1451.8

# CRP!A126 row
$ws.Range("H126").Value = 44016.582
$ws.Range("I126").Value = 73209.57000000001
$ws.Range("J126").Value = 3146.4
$ws.Range("K126").Value = 219628.71
$ws.Range("L126").Value = 9439.200000000001
$ws.Range("M126").Value = -217158.71
$ws.Range("N126").Value = -14379.2

# CRP!A132 row
$ws.Range("H132").Value = 2119.2354
$ws.Range("I132").Value = 1108.7931
$ws.Range("J132").Value = 3451.182
$ws.Range("K132").Value = 3326.379300000001
$ws.Range("L132").Value = 10353.546
$ws.Range("M132").Value = -796.3793000000005
$ws.Range("N132").Value = -15413.546

# CRP!A134 row
$ws.Range("H134").Value = 1918.75
$ws.Range("I134").Value = 1122.7059
$ws.Range("J134").Value = 2820.9333
$ws.Range("K134").Value = 3368.1177
$ws.Range("L134").Value = 8462.7999
$ws.Range("M134").Value = -833.1176999999998
$ws.Range("N134").Value = -13532.7999

# CRP!A135 row
$ws.Range("H135").Value = 39567.855
$ws.Range("J135").Value = 28000
$ws.Range("L135").Value = 28000
$ws.Range("N135").Value = -38140

# CUL!A68 row
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 472.8
$ws.Range("I68").Value = 426.66666
$ws.Range("J68").Value = 492.57144
$ws.Range("K68").Value = 1279.99998
$ws.Range("L68").Value = 1477.71432
$ws.Range("M68").Value = -468.9999800000001
$ws.Range("N68").Value = -3099.71432

# CUL!A71 row
$ws.Range("H71").Value = 472.8
$ws.Range("I71").Value = 426.66666
$ws.Range("J71").Value = 492.57144
$ws.Range("K71").Value = 3839.99994
$ws.Range("L71").Value = 4433.14296
$ws.Range("M71").Value = 216.0000600000003
$ws.Range("N71").Value = -12545.14296

# CUL!A109 row
$ws.Range("H109").Value = 2866.6667
$ws.Range("I109").Value = 1333.3334
$ws.Range("J109").Value = 3633.3333
$ws.Range("K109").Value = 4000.0002
$ws.Range("L109").Value = 10899.9999
$ws.Range("M109").Value = -2960.0002
$ws.Range("N109").Value = -12979.9999

# CUL!A112 row
$ws.Range("H112").Value = 4551
$ws.Range("J112").Value = 4551
$ws.Range("L112").Value = 13653
$ws.Range("N112").Value = -15869

# GSM!A102 row
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1133.3334
$ws.Range("I102").Value = 1112.5312
$ws.Range("J102").Value = 1228.4286
$ws.Range("K102").Value = 1112.5312
$ws.Range("L102").Value = 1228.4286
$ws.Range("M102").Value = 509.4688000000001
$ws.Range("N102").Value = -4472.4286

# GSM!A113 row
$ws.Range("H113").Value = 1139
$ws.Range("I113").Value = 1139
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1139
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("M113").Value = 1031

# GSM!A126 row
$ws.Range("H126").Value = 2753.625
$ws.Range("I126").Value = 2618.1
$ws.Range("J126").Value = 2979.5
$ws.Range("K126").Value = 7854.299999999999
$ws.Range("L126").Value = 8938.5
$ws.Range("M126").Value = -5384.299999999999
$ws.Range("N126").Value = -13878.5

# GSM!A132 row
$ws.Range("H132").Value = 3373.0833
$ws.Range("I132").Value = 3427.4688
$ws.Range("J132").Value = 3264.3125
$ws.Range("K132").Value = 10282.4064
$ws.Range("L132").Value = 9792.9375
$ws.Range("M132").Value = -7752.4064
$ws.Range("N132").Value = -14852.9375

# GSM!A133 row
$ws.Range("H133").Value = 24508
$ws.Range("J133").Value = 24508
$ws.Range("L133").Value = 24508
$ws.Range("N133").Value = -34628

# LTW!A111 row
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 23500
$ws.Range("J111").Value = 23500
$ws.Range("L111").Value = 23500
$ws.Range("N111").Value = -31680

# WVR!A107 row
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 369.6316
$ws.Range("I107").Value = 369.13333
$ws.Range("J107").Value = 371.5
$ws.Range("K107").Value = 1107.39999
$ws.Range("L107").Value = 1114.5
$ws.Range("M107").Value = 812.6000100000001
$ws.Range("N107").Value = -4954.5

# WVR!A126 row
$ws.Range("H126").Value = 1506.7097
$ws.Range("I126").Value = 1145.2916
$ws.Range("J126").Value = 2745.8572
$ws.Range("K126").Value = 3435.8748
$ws.Range("L126").Value = 8237.571599999999
$ws.Range("M126").Value = -965.8748000000001
$ws.Range("N126").Value = -13177.5716

# WVR!A132 row
$ws.Range("H132").Value = 16472.426
$ws.Range("I132").Value = 25483.023
$ws.Range("J132").Value = 2789.6667
$ws.Range("K132").Value = 76449.069
$ws.Range("L132").Value = 8369.000100000001
$ws.Range("M132").Value = -73919.069
$ws.Range("N132").Value = -13429.0001

# WVR!A136 row
$ws.Range("H136").Value = 19232372
$ws.Range("I136").Value = 38462830
$ws.Range("J136").Value = 1911.3462
$ws.Range("K136").Value = 115388490
$ws.Range("L136").Value = 5734.0386
$ws.Range("M136").Value = -115385940
$ws.Range("N136").Value = -10834.0386
